$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Data Type" and "Number of Unique Values" columns (old columns B & C),
# leaving "Column" (A) and "Description" (old D, now B)
$ws.Range("B1:C1").EntireColumn.Delete()

# Resize the table to the new 2-column extent and resync the header names
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B7"))
$ws.Range("A1").Value = "Column"
$ws.Range("B1").Value = "Description"

# Update a couple of description texts that changed content
$ws.Range("B3").Value = "Annual; Ranges from 2000 to 2024"
$ws.Range("B7").Value = "Population estimates"

# Row 7 no longer needs its taller custom height
$ws.Rows.Item(7).AutoFit()

# Restore the view's selection (frozen pane auto-clamps to column B already)
$ws.Range("B14").Select() | Out-Null
